$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 949.1053000000001
$ws.Cells.Item(19, 9).Value = 906.3333
$ws.Cells.Item(19, 11).Value = 906.3333
$ws.Cells.Item(19, 13).Value = -731.3333
$ws.Cells.Item(28, 8).Value = 918.2222
$ws.Cells.Item(28, 10).Value = 3661.6667
$ws.Cells.Item(28, 12).Value = 3661.6667
$ws.Cells.Item(28, 14).Value = -4631.6667
$ws.Cells.Item(111, 8).Value = 810.1539
$ws.Cells.Item(111, 10).Value = 2110
$ws.Cells.Item(111, 12).Value = 6330
$ws.Cells.Item(111, 14).Value = -12464
$ws.Cells.Item(112, 8).Value = 1851.6471
$ws.Cells.Item(112, 10).Value = 1961.7333
$ws.Cells.Item(112, 12).Value = 5885.199900000001
$ws.Cells.Item(112, 14).Value = -8101.199900000001
$ws.Cells.Item(116, 9).Value = 6692.6665
$ws.Cells.Item(116, 10).Value = 9301
$ws.Cells.Item(116, 11).Value = 6692.6665
$ws.Cells.Item(116, 12).Value = 9301
$ws.Cells.Item(116, 13).Value = -3250.6665
$ws.Cells.Item(116, 14).Value = -16185
$ws.Cells.Item(129, 8).Value = 1060.0667
$ws.Cells.Item(129, 9).Value = 477.45456
$ws.Cells.Item(129, 11).Value = 1432.36368
$ws.Cells.Item(129, 13).Value = 3567.63632
$ws.Cells.Item(138, 8).Value = 3350.2207
$ws.Cells.Item(138, 9).Value = 1194.75
$ws.Cells.Item(138, 11).Value = 3584.25
$ws.Cells.Item(138, 13).Value = 1555.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4859.644
$ws.Cells.Item(32, 9).Value = 3570.8086
$ws.Cells.Item(32, 11).Value = 3570.8086
$ws.Cells.Item(32, 13).Value = -3283.8086
$ws.Cells.Item(102, 8).Value = 1511.25
$ws.Cells.Item(102, 9).Value = 1198.8334
$ws.Cells.Item(102, 10).Value = 2448.5
$ws.Cells.Item(102, 11).Value = 1198.8334
$ws.Cells.Item(102, 12).Value = 2448.5
$ws.Cells.Item(102, 13).Value = 423.1666
$ws.Cells.Item(102, 14).Value = -5692.5
$ws.Cells.Item(122, 8).Value = 3991.8333
$ws.Cells.Item(122, 9).Value = 3578.6667
$ws.Cells.Item(122, 11).Value = 10736.0001
$ws.Cells.Item(122, 13).Value = -8286.000100000001
$ws.Cells.Item(132, 8).Value = 1945.1833
$ws.Cells.Item(132, 9).Value = 1561.2
$ws.Cells.Item(132, 11).Value = 4683.6
$ws.Cells.Item(132, 13).Value = -2153.6
$ws.Cells.Item(138, 8).Value = 79996.336
$ws.Cells.Item(138, 10).Value = 79996.336
$ws.Cells.Item(138, 12).Value = 79996.336
$ws.Cells.Item(138, 14).Value = -90276.336
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 857.8
$ws.Cells.Item(64, 10).Value = 897.5
$ws.Cells.Item(64, 12).Value = 897.5
$ws.Cells.Item(64, 14).Value = -1347.5
$ws.Cells.Item(67, 8).Value = 857.8
$ws.Cells.Item(67, 10).Value = 897.5
$ws.Cells.Item(67, 12).Value = 897.5
$ws.Cells.Item(67, 14).Value = -2457.5
$ws.Cells.Item(99, 8).Value = 2992.3333
$ws.Cells.Item(99, 9).Value = 2561.5715
$ws.Cells.Item(99, 11).Value = 2561.5715
$ws.Cells.Item(99, 13).Value = -1063.5715
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(103, 8).Value = 30535.615
$ws.Cells.Item(103, 9).Value = 28913.584
$ws.Cells.Item(103, 11).Value = 28913.584
$ws.Cells.Item(103, 13).Value = -27741.584
$ws.Cells.Item(122, 8).Value = 5134.3
$ws.Cells.Item(122, 9).Value = 2019.909
$ws.Cells.Item(122, 10).Value = 8940.777
$ws.Cells.Item(122, 11).Value = 6059.727000000001
$ws.Cells.Item(122, 12).Value = 26822.331
$ws.Cells.Item(122, 13).Value = -3609.727000000001
$ws.Cells.Item(122, 14).Value = -31722.331
$ws.Cells.Item(132, 8).Value = 4566.12
$ws.Cells.Item(132, 9).Value = 3875.6086
$ws.Cells.Item(132, 11).Value = 11626.8258
$ws.Cells.Item(132, 13).Value = -9096.825800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 5953.75
$ws.Cells.Item(88, 9).Value = 2933
$ws.Cells.Item(88, 10).Value = 15016
$ws.Cells.Item(88, 11).Value = 8799
$ws.Cells.Item(88, 12).Value = 45048
$ws.Cells.Item(88, 13).Value = -8371
$ws.Cells.Item(88, 14).Value = -45904
$ws.Cells.Item(91, 8).Value = 5953.75
$ws.Cells.Item(91, 9).Value = 2933
$ws.Cells.Item(91, 10).Value = 15016
$ws.Cells.Item(91, 11).Value = 8799
$ws.Cells.Item(91, 12).Value = 45048
$ws.Cells.Item(91, 13).Value = -7317
$ws.Cells.Item(91, 14).Value = -48012
$ws.Cells.Item(107, 8).Value = 440514.2
$ws.Cells.Item(107, 9).Value = 375.8
$ws.Cells.Item(107, 10).Value = 1202292.1
$ws.Cells.Item(107, 11).Value = 1127.4
$ws.Cells.Item(107, 12).Value = 3606876.3
$ws.Cells.Item(107, 13).Value = 792.5999999999999
$ws.Cells.Item(107, 14).Value = -3610716.3
$ws.Cells.Item(117, 8).Value = 3834.7646
$ws.Cells.Item(117, 10).Value = 4107.846
$ws.Cells.Item(117, 12).Value = 12323.538
$ws.Cells.Item(117, 14).Value = -19207.538
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 570.13336
$ws.Cells.Item(2, 9).Value = 49.090908
$ws.Cells.Item(2, 11).Value = 49.090908
$ws.Cells.Item(2, 13).Value = 63.909092
$ws.Cells.Item(54, 8).Value = 8299.799999999999
$ws.Cells.Item(54, 9).Value = 6249.5
$ws.Cells.Item(54, 11).Value = 6249.5
$ws.Cells.Item(54, 13).Value = -5859.5
$ws.Cells.Item(102, 8).Value = 3474.7827
$ws.Cells.Item(102, 9).Value = 1733.8667
$ws.Cells.Item(102, 11).Value = 1733.8667
$ws.Cells.Item(102, 13).Value = -111.8667
$ws.Cells.Item(132, 8).Value = 3382.6458
$ws.Cells.Item(132, 9).Value = 2373.8235
$ws.Cells.Item(132, 11).Value = 7121.470499999999
$ws.Cells.Item(132, 13).Value = -4591.470499999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 19999.924
$ws.Cells.Item(2, 9).Value = 19222.111
$ws.Cells.Item(2, 11).Value = 19222.111
$ws.Cells.Item(2, 13).Value = -19110.111
$ws.Cells.Item(16, 8).Value = 961.6
$ws.Cells.Item(16, 9).Value = 485.94736
$ws.Cells.Item(16, 11).Value = 485.94736
$ws.Cells.Item(16, 13).Value = -315.94736
$ws.Cells.Item(61, 8).Value = 3409.1035
$ws.Cells.Item(61, 9).Value = 2254.037
$ws.Cells.Item(61, 11).Value = 2254.037
$ws.Cells.Item(61, 13).Value = -2052.037
$ws.Cells.Item(100, 8).Value = 5228.231
$ws.Cells.Item(100, 9).Value = 4220.5
$ws.Cells.Item(100, 11).Value = 4220.5
$ws.Cells.Item(100, 13).Value = -3679.5
$ws.Cells.Item(113, 8).Value = 3409.1035
$ws.Cells.Item(113, 9).Value = 2254.037
$ws.Cells.Item(113, 11).Value = 2254.037
$ws.Cells.Item(113, 13).Value = -84.03699999999981
$ws.Cells.Item(122, 8).Value = 8306.299999999999
$ws.Cells.Item(122, 9).Value = 7294.143
$ws.Cells.Item(122, 11).Value = 21882.429
$ws.Cells.Item(122, 13).Value = -19432.429
$ws.Cells.Item(132, 8).Value = 3577.8909
$ws.Cells.Item(132, 9).Value = 2733.861
$ws.Cells.Item(132, 11).Value = 8201.582999999999
$ws.Cells.Item(132, 13).Value = -5671.582999999999
$ws.Cells.Item(136, 8).Value = 4431.0967
$ws.Cells.Item(136, 9).Value = 2639.0454
$ws.Cells.Item(136, 11).Value = 7917.1362
$ws.Cells.Item(136, 13).Value = -5367.1362
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2246.5588
$ws.Cells.Item(122, 9).Value = 1675.3448
$ws.Cells.Item(122, 10).Value = 5559.6
$ws.Cells.Item(122, 11).Value = 5026.0344
$ws.Cells.Item(122, 12).Value = 16678.8
$ws.Cells.Item(122, 13).Value = -2576.0344
$ws.Cells.Item(122, 14).Value = -21578.8
